$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.09665433333333333
$ws.Range("H2").Value = 0.289963
$ws.Range("I2").Value = 0.0006230336790718351
$ws.Range("J2").Value = 0.0006230336790718351
$ws.Range("M2").Value = 110.642708
$ws.Range("N2").Value = 331.928124
$ws.Range("O2").Value = 0.5476418925386564
$ws.Range("P2").Value = 0.5476418925386564
$ws.Range("Q2").Value = 10.69409717993466
$ws.Range("R2").Value = 96.24687461941198
$ws.Range("S2").Value = 0.0003411993431222216
$ws.Range("T2").Value = 0.0003411993431222216

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.09665433333333333
$ws.Range("H3").Value = 0.289963
$ws.Range("I3").Value = 0.0006230336790718351
$ws.Range("J3").Value = 0.0006230336790718351
$ws.Range("O3").Value = 0.3151072754333865
$ws.Range("P3").Value = 0.3151072754333865
$ws.Range("Q3").Value = 6.153268899806111
$ws.Range("R3").Value = 55.379420098255
$ws.Range("S3").Value = 0.0001963224451155649
$ws.Range("T3").Value = 0.0001963224451155649

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.09665433333333333
$ws.Range("H4").Value = 0.289963
$ws.Range("I4").Value = 0.0006230336790718351
$ws.Range("J4").Value = 0.0006230336790718351
$ws.Range("M4").Value = 27.72944133333333
$ws.Range("N4").Value = 83.188324
$ws.Range("O4").Value = 0.1372508320279571
$ws.Range("P4").Value = 0.1372508320279571
$ws.Range("Q4").Value = 2.680170665779111
$ws.Range("R4").Value = 24.121535992012
$ws.Range("S4").Value = 0.000085511890834048578913760552
$ws.Range("T4").Value = 0.000085511890834048578913760552

$ws.Range("G5").Value = 154.8642143333334
$ws.Range("H5").Value = 464.5926430000001
$ws.Range("I5").Value = 0.9982544794956518
$ws.Range("J5").Value = 0.9982544794956519
$ws.Range("M5").Value = 110.642708
$ws.Range("N5").Value = 331.928124
$ws.Range("O5").Value = 0.5476418925386564
$ws.Range("P5").Value = 0.5476418925386564
$ws.Range("Q5").Value = 17134.59604613242
$ws.Range("R5").Value = 154211.3644151917
$ws.Range("S5").Value = 0.5466859723861901
$ws.Range("T5").Value = 0.5466859723861902

$ws.Range("G6").Value = 154.8642143333334
$ws.Range("H6").Value = 464.5926430000001
$ws.Range("I6").Value = 0.9982544794956518
$ws.Range("J6").Value = 0.9982544794956519
$ws.Range("O6").Value = 0.3151072754333865
$ws.Range("P6").Value = 0.3151072754333865
$ws.Range("Q6").Value = 9859.062919236676
$ws.Range("R6").Value = 88731.56627313007
$ws.Range("S6").Value = 0.3145572492230482
$ws.Range("T6").Value = 0.3145572492230483

$ws.Range("G7").Value = 154.8642143333334
$ws.Range("H7").Value = 464.5926430000001
$ws.Range("I7").Value = 0.9982544794956518
$ws.Range("J7").Value = 0.9982544794956519
$ws.Range("M7").Value = 27.72944133333333
$ws.Range("N7").Value = 83.188324
$ws.Range("O7").Value = 0.1372508320279571
$ws.Range("P7").Value = 0.1372508320279571
$ws.Range("Q7").Value = 4294.298145988926
$ws.Range("R7").Value = 38648.68331390034
$ws.Range("S7").Value = 0.1370112578864135
$ws.Range("T7").Value = 0.1370112578864135

$ws.Range("G8").Value = 0.174137
$ws.Range("H8").Value = 0.522411
$ws.Range("I8").Value = 0.001122486825276316
$ws.Range("J8").Value = 0.001122486825276316
$ws.Range("M8").Value = 110.642708
$ws.Range("N8").Value = 331.928124
$ws.Range("O8").Value = 0.5476418925386564
$ws.Range("P8").Value = 0.5476418925386564
$ws.Range("Q8").Value = 19.266989242996
$ws.Range("R8").Value = 173.402903186964
$ws.Range("S8").Value = 0.0006147208093440296
$ws.Range("T8").Value = 0.0006147208093440299

$ws.Range("G9").Value = 0.174137
$ws.Range("H9").Value = 0.522411
$ws.Range("I9").Value = 0.001122486825276316
$ws.Range("J9").Value = 0.001122486825276316
$ws.Range("O9").Value = 0.3151072754333865
$ws.Range("P9").Value = 0.3151072754333865
$ws.Range("Q9").Value = 11.08601911008167
$ws.Range("R9").Value = 99.774171990735
$ws.Range("S9").Value = 0.0003537037652226917
$ws.Range("T9").Value = 0.0003537037652226917

$ws.Range("G10").Value = 0.174137
$ws.Range("H10").Value = 0.522411
$ws.Range("I10").Value = 0.001122486825276316
$ws.Range("J10").Value = 0.001122486825276316
$ws.Range("M10").Value = 27.72944133333333
$ws.Range("N10").Value = 83.188324
$ws.Range("O10").Value = 0.1372508320279571
$ws.Range("P10").Value = 0.1372508320279571
$ws.Range("Q10").Value = 4.828721725462666
$ws.Range("R10").Value = 43.458495529164
$ws.Range("S10").Value = 0.0001540622507095945
$ws.Range("T10").Value = 0.0001540622507095945

